# Auto-assembled PowerShell Excel COM-interop script
# Applies the commit diff: adds 'Player Info' and 'ODI Batting Extra' sheets,
# renames MATCH_CARD_LINK -> MATCH_CODE columns (URL -> numeric code) on the
# existing 'ODI Batting' / 'ODI Bowling' sheets, and drops stray empty INNING_NUMBER cells.

$wb = $excel.ActiveWorkbook

# --- 1. Insert "Player Info" as the new first sheet -----------------------
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$playerInfo = $wb.Worksheets.Add($battingSheet)
$playerInfo.Name = "Player Info"

$piHeaders = @("ID","NAME","BATTING_HAND","BOWL_STYLE")
for ($c = 1; $c -le $piHeaders.Length; $c++) {
    $cell = $playerInfo.Cells.Item(1, $c)
    $cell.Value = $piHeaders[$c - 1]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

$playerInfo.Range("A2").NumberFormat = "@"
$piRow = @("4310","Axar Rajeshbhai Patel","Left Handed","Left Arm Orthodox")
for ($c = 1; $c -le $piRow.Length; $c++) {
    $playerInfo.Cells.Item(2, $c).Value = $piRow[$c - 1]
}

# --- 2. "ODI Batting": MATCH_CARD_LINK -> MATCH_CODE, URL -> bare code -----
$battingSheet.Cells.Item(1, 4).Value = "MATCH_CODE"

$battingSheet.Range("D2:D52").NumberFormat = "@"
$battingMatchCodes = @{
    2 = "3643"
    3 = "3644"
    4 = "3645"
    5 = "3683"
    6 = "3688"
    7 = "3689"
    8 = "3692"
    9 = "3693"
    10 = "3696"
    11 = "3731"
    12 = "3735"
    13 = "3741"
    14 = "3744"
    15 = "3810"
    16 = "3811"
    17 = "3812"
    18 = "3815"
    19 = "3817"
    20 = "3844"
    21 = "3845"
    22 = "3848"
    23 = "3851"
    24 = "3895"
    25 = "3897"
    26 = "3899"
    27 = "3951"
    28 = "3952"
    29 = "3953"
    30 = "3954"
    31 = "3955"
    32 = "4061"
    33 = "4062"
    34 = "4063"
    35 = "4064"
    36 = "4074"
    37 = "4076"
    38 = "4088"
    39 = "4089"
    40 = "4621"
    41 = "4623"
    42 = "4624"
    43 = "4637"
    44 = "4640"
    45 = "4643"
    46 = "4682"
    47 = "4685"
    48 = "4687"
    49 = "4689"
    50 = "4691"
    51 = "4728"
    52 = "4732"
}
foreach ($r in $battingMatchCodes.Keys) {
    $battingSheet.Cells.Item($r, 4).Value = $battingMatchCodes[$r]
}

# Rows whose INNING_NUMBER (col B) cell was a stray empty cell -> drop it entirely
$battingEmptyBRows = @(2,5,7,8,9,13,19,24,25,26,27,29,32,34,35,37,38,39,42,43)
foreach ($r in $battingEmptyBRows) {
    $battingSheet.Cells.Item($r, 2).ClearContents()
}

# --- 3. "ODI Bowling": MATCH_CARD_LINK -> MATCH_CODE, URL -> bare code ------
$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")
$bowlingSheet.Cells.Item(1, 2).Value = "MATCH_CODE"

$bowlingSheet.Range("B2:B48").NumberFormat = "@"
$bowlingMatchCodes = @{
    2 = "3643"
    3 = "3683"
    4 = "3688"
    5 = "3689"
    6 = "3692"
    7 = "3693"
    8 = "3696"
    9 = "3731"
    10 = "3735"
    11 = "3744"
    12 = "3810"
    13 = "3811"
    14 = "3812"
    15 = "3815"
    16 = "3817"
    17 = "3844"
    18 = "3845"
    19 = "3848"
    20 = "3851"
    21 = "3895"
    22 = "3897"
    23 = "3899"
    24 = "3951"
    25 = "3952"
    26 = "3953"
    27 = "3954"
    28 = "3955"
    29 = "4061"
    30 = "4062"
    31 = "4063"
    32 = "4064"
    33 = "4074"
    34 = "4076"
    35 = "4088"
    36 = "4089"
    37 = "4621"
    38 = "4623"
    39 = "4624"
    40 = "4637"
    41 = "4640"
    42 = "4643"
    43 = "4682"
    44 = "4685"
    45 = "4687"
    46 = "4689"
    47 = "4728"
    48 = "4732"
}
foreach ($r in $bowlingMatchCodes.Keys) {
    $bowlingSheet.Cells.Item($r, 2).Value = $bowlingMatchCodes[$r]
}

# --- 4. Append "ODI Batting Extra" as the new last sheet --------------------
$extra = $wb.Worksheets.Add($null, $bowlingSheet)
$extra.Name = "ODI Batting Extra"

$extraHeaders = @("MATCH_CODE","BATTING_POSITION","NUM_4","NUM_6","PERCENT_RUNS_OF_TOTAL","MAN_OF_MATCH")
for ($c = 1; $c -le $extraHeaders.Length; $c++) {
    $cell = $extra.Cells.Item(1, $c)
    $cell.Value = $extraHeaders[$c - 1]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

$extra.Range("A2:A21").NumberFormat = "@"
$extra.Range("C2:E21").NumberFormat = "@"

$extraData = @(
    @("4062", "NONE", $null, $null, $null, "NO"),
    @("4063", "NONE", $null, $null, $null, "NO"),
    @("4064", 9, $null, $null, $null, "NO"),
    @("4074", "NONE", $null, $null, $null, "NO"),
    @("4076", 9, $null, $null, $null, "NO"),
    @("4088", "NONE", $null, $null, $null, "NO"),
    @("4089", 8, $null, $null, $null, "NO"),
    @("4621", 7, "1", "1", "6.82%", "NO"),
    @("4623", "NONE", $null, $null, $null, "NO"),
    @("4624", 7, $null, $null, $null, "NO"),
    @("4637", 7, $null, $null, $null, "NO"),
    @("4640", 7, "1", "0", "3.59%", "NO"),
    @("4643", "NONE", $null, $null, $null, "NO"),
    @("4682", 6, "2", "3", "21.05%", "NO"),
    @("4685", "NONE", $null, $null, $null, "NO"),
    @("4687", "NONE", $null, $null, $null, "NO"),
    @("4689", 7, "1", "1", "9.59%", "NO"),
    @("4691", "NONE", $null, $null, $null, "NO"),
    @("4728", 8, "1", "2", "24.79%", "NO"),
    @("4732", 5, "0", "0", "0.81%", "NO"),
)

$r = 2
foreach ($row in $extraData) {
    $extra.Cells.Item($r, 1).Value = $row[0]
    if ($row[1] -ne "NONE") {
        $extra.Cells.Item($r, 2).Value = [int]$row[1]
    }
    if ($row[2] -ne $null) { $extra.Cells.Item($r, 3).Value = $row[2] }
    if ($row[3] -ne $null) { $extra.Cells.Item($r, 4).Value = $row[3] }
    if ($row[4] -ne $null) { $extra.Cells.Item($r, 5).Value = $row[4] }
    $extra.Cells.Item($r, 6).Value = $row[5]
    $r++
}

"done"
